$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "2022_08_06"

$ws.Range("D2").Select()
